$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = 'cum_Govt'
$ws.Activate()

$ws.Range("A1").Value = 'Country'
$ws.Range("B1").Value = 'cum_Govt -> PctUrb chi2'
$ws.Range("C1").Value = 'Prob > chi2'
$ws.Range("D1").Value = 'PctUrb -> cum_Govt chi2'
$ws.Range("E1").Value = 'Prob > chi2'

$ws.Range("A2").Value = 'France'
$ws.Range("B2").Value = [double]"116.76215944202019"
$ws.Range("C2").Value = [double]"4.4199690634033109E-26"
$ws.Range("D2").Value = [double]"0.34583818699846536"
$ws.Range("E2").Value = [double]"0.84120567109876365"
$ws.Range("A3").Value = 'Azerbaijan'
$ws.Range("B3").Value = [double]"102.82456765729553"
$ws.Range("C3").Value = [double]"4.6981711229677782E-23"
$ws.Range("D3").Value = [double]"0.91289110816188934"
$ws.Range("E3").Value = [double]"0.63353150165939476"
$ws.Range("A4").Value = 'Vietnam'
$ws.Range("B4").Value = [double]"72.848306739473145"
$ws.Range("C4").Value = [double]"1.5177183448304706E-16"
$ws.Range("D4").Value = [double]"3.7520501354159883"
$ws.Range("E4").Value = [double]"0.1531978481629567"
$ws.Range("A5").Value = 'Mongolia'
$ws.Range("B5").Value = [double]"57.667767402888558"
$ws.Range("C5").Value = [double]"3.0033325059937359E-13"
$ws.Range("D5").Value = [double]"1.3213564060092864"
$ws.Range("E5").Value = [double]"0.51650092320211916"
$ws.Range("A6").Value = 'Thailand'
$ws.Range("B6").Value = [double]"34.854201293550467"
$ws.Range("C6").Value = [double]"2.7008866207920427E-8"
$ws.Range("D6").Value = [double]"4.5557574889317376"
$ws.Range("E6").Value = [double]"0.10250140794329474"
$ws.Range("A7").Value = 'Panama'
$ws.Range("B7").Value = [double]"31.494989971924269"
$ws.Range("C7").Value = [double]"1.4486044793926115E-7"
$ws.Range("D7").Value = [double]"46.942785548923027"
$ws.Range("E7").Value = [double]"6.4047714296395926E-11"
$ws.Range("A8").Value = 'Bangladesh'
$ws.Range("B8").Value = [double]"27.582335438960968"
$ws.Range("C8").Value = [double]"1.0246415442336814E-6"
$ws.Range("D8").Value = [double]"28.815183225303155"
$ws.Range("E8").Value = [double]"5.5317489884523327E-7"
$ws.Range("A9").Value = 'El Salvador'
$ws.Range("B9").Value = [double]"27.38097097847583"
$ws.Range("C9").Value = [double]"1.1331768600113821E-6"
$ws.Range("D9").Value = [double]"1.7694736027912539"
$ws.Range("E9").Value = [double]"0.41282281336979465"
$ws.Range("A10").Value = 'Zambia'
$ws.Range("B10").Value = [double]"27.170764617664418"
$ws.Range("C10").Value = [double]"1.258761434639437E-6"
$ws.Range("D10").Value = [double]"2.8833481343333993"
$ws.Range("E10").Value = [double]"0.23653145751037502"
$ws.Range("A11").Value = 'Nigeria'
$ws.Range("B11").Value = [double]"26.869276699712142"
$ws.Range("C11").Value = [double]"1.4635605622856722E-6"
$ws.Range("D11").Value = [double]"0.90522315802600506"
$ws.Range("E11").Value = [double]"0.63596510785321114"
$ws.Range("A12").Value = 'Indonesia'
$ws.Range("B12").Value = [double]"26.144564524249972"
$ws.Range("C12").Value = [double]"2.1027127479665662E-6"
$ws.Range("D12").Value = [double]"13.807594086856909"
$ws.Range("E12").Value = [double]"1.0039660797142953E-3"
$ws.Range("A13").Value = 'Burkina Faso'
$ws.Range("B13").Value = [double]"25.673171659319969"
$ws.Range("C13").Value = [double]"2.6615927673305031E-6"
$ws.Range("D13").Value = [double]"2.3965071822470776"
$ws.Range("E13").Value = [double]"0.30172067973765798"
$ws.Range("A14").Value = 'Hungary'
$ws.Range("B14").Value = [double]"24.692620705878198"
$ws.Range("C14").Value = [double]"4.3457580353237854E-6"
$ws.Range("D14").Value = [double]"7.2792452694655054"
$ws.Range("E14").Value = [double]"2.626225255794614E-2"
$ws.Range("A15").Value = 'Greece'
$ws.Range("B15").Value = [double]"23.58713825133097"
$ws.Range("C15").Value = [double]"7.5529742969390981E-6"
$ws.Range("D15").Value = [double]"17.265878145246745"
$ws.Range("E15").Value = [double]"1.7814030770811067E-4"
$ws.Range("A16").Value = 'Malawi'
$ws.Range("B16").Value = [double]"23.236393449968382"
$ws.Range("C16").Value = [double]"9.0008033366454049E-6"
$ws.Range("D16").Value = [double]"8.1998481286077318"
$ws.Range("E16").Value = [double]"1.6573933907186821E-2"
$ws.Range("A17").Value = 'Brazil'
$ws.Range("B17").Value = [double]"22.671577096590852"
$ws.Range("C17").Value = [double]"1.1937945323158705E-5"
$ws.Range("D17").Value = [double]"0.98842674022958565"
$ws.Range("E17").Value = [double]"0.61005060263575062"
$ws.Range("A18").Value = 'Germany'
$ws.Range("B18").Value = [double]"20.793777586986835"
$ws.Range("C18").Value = [double]"3.0527312187638162E-5"
$ws.Range("D18").Value = [double]"10.959870633290805"
$ws.Range("E18").Value = [double]"4.1695993928576895E-3"
$ws.Range("A19").Value = 'United Kingdom'
$ws.Range("B19").Value = [double]"20.702486555494716"
$ws.Range("C19").Value = [double]"3.195303857009905E-5"
$ws.Range("D19").Value = [double]"1.6948473049205717"
$ws.Range("E19").Value = [double]"0.42851752107360386"
$ws.Range("A20").Value = 'Australia'
$ws.Range("B20").Value = [double]"18.914386167846828"
$ws.Range("C20").Value = [double]"7.812557546964244E-5"
$ws.Range("D20").Value = [double]"24.372611479507981"
$ws.Range("E20").Value = [double]"5.0998178235018482E-6"
$ws.Range("A21").Value = 'Belgium'
$ws.Range("B21").Value = [double]"18.718253358557757"
$ws.Range("C21").Value = [double]"8.6175324560002755E-5"
$ws.Range("D21").Value = [double]"63.269179414885265"
$ws.Range("E21").Value = [double]"1.803034991941925E-15"
$ws.Range("A22").Value = 'Haiti'
$ws.Range("B22").Value = [double]"18.638809117306792"
$ws.Range("C22").Value = [double]"8.9667286119957365E-5"
$ws.Range("D22").Value = [double]"5.1038460180008229"
$ws.Range("E22").Value = [double]"7.793165853374763E-2"
$ws.Range("A23").Value = 'Namibia'
$ws.Range("B23").Value = [double]"18.260168355724602"
$ws.Range("C23").Value = [double]"1.0835646349993803E-4"
$ws.Range("D23").Value = [double]"3.5363892822702976"
$ws.Range("E23").Value = [double]"0.17064077874879652"
$ws.Range("A24").Value = 'Switzerland'
$ws.Range("B24").Value = [double]"16.94864016818542"
$ws.Range("C24").Value = [double]"2.0876108703626447E-4"
$ws.Range("D24").Value = [double]"0.10482523600859758"
$ws.Range("E24").Value = [double]"0.94893723746265257"
$ws.Range("A25").Value = 'Iraq'
$ws.Range("B25").Value = [double]"14.271995634081321"
$ws.Range("C25").Value = [double]"7.9593119281712603E-4"
$ws.Range("D25").Value = [double]"0.26657350848911582"
$ws.Range("E25").Value = [double]"0.87521408476807683"
$ws.Range("A26").Value = 'Venezuela, RB'
$ws.Range("B26").Value = [double]"13.583608427612489"
$ws.Range("C26").Value = [double]"1.1229409198723289E-3"
$ws.Range("D26").Value = [double]"0.63779534152252149"
$ws.Range("E26").Value = [double]"0.72694993373374295"
$ws.Range("A27").Value = 'Botswana'
$ws.Range("B27").Value = [double]"12.503677085764139"
$ws.Range("C27").Value = [double]"1.9269081742168443E-3"
$ws.Range("D27").Value = [double]"30.53022997110736"
$ws.Range("E27").Value = [double]"2.346630959830227E-7"
$ws.Range("A28").Value = 'Tajikistan'
$ws.Range("B28").Value = [double]"10.599318579386818"
$ws.Range("C28").Value = [double]"4.9932948841544022E-3"
$ws.Range("D28").Value = [double]"7.6034986398118489"
$ws.Range("E28").Value = [double]"2.2331672428321174E-2"
$ws.Range("A29").Value = 'Togo'
$ws.Range("B29").Value = [double]"9.8994234423026199"
$ws.Range("C29").Value = [double]"7.0854512203835408E-3"
$ws.Range("D29").Value = [double]"20065.354695057918"
$ws.Range("E29").Value = [double]"0"
$ws.Range("A30").Value = 'Kyrgyz Republic'
$ws.Range("B30").Value = [double]"9.1280471924310405"
$ws.Range("C30").Value = [double]"1.0420048416031827E-2"
$ws.Range("D30").Value = [double]"6.5352219702027394"
$ws.Range("E30").Value = [double]"3.809733354091456E-2"
$ws.Range("A31").Value = 'Ghana'
$ws.Range("B31").Value = [double]"9.117718324576316"
$ws.Range("C31").Value = [double]"1.0474001265643495E-2"
$ws.Range("D31").Value = [double]"5.6312565516643787"
$ws.Range("E31").Value = [double]"5.9867093873728426E-2"
$ws.Range("A32").Value = 'Papua New Guinea'
$ws.Range("B32").Value = [double]"5.9535474616641295"
$ws.Range("C32").Value = [double]"5.0956969849810914E-2"
$ws.Range("D32").Value = [double]"5.1872356792345116"
$ws.Range("E32").Value = [double]"7.4749119980463088E-2"
$ws.Range("A33").Value = 'Cote d''Ivoire'
$ws.Range("B33").Value = [double]"5.569942726813542"
$ws.Range("C33").Value = [double]"6.1730856767082627E-2"
$ws.Range("D33").Value = [double]"6.1803039765349244"
$ws.Range("E33").Value = [double]"4.5495039166936271E-2"
$ws.Range("A34").Value = 'Iceland'
$ws.Range("B34").Value = [double]"2.3169668414621327"
$ws.Range("C34").Value = [double]"0.31396196821904138"
$ws.Range("D34").Value = [double]"3.2516663199026326"
$ws.Range("E34").Value = [double]"0.1967476846071525"
$ws.Range("A35").Value = 'Kenya'
$ws.Range("B35").Value = [double]"0.40661535478068345"
$ws.Range("C35").Value = [double]"0.81602712969917668"
$ws.Range("D35").Value = [double]"12.184681559945659"
$ws.Range("E35").Value = [double]"2.2601122925284003E-3"

$ws.Range("B2:E35").NumberFormat = "0.000"

$ws.Columns.Item(1).ColumnWidth = 18
$ws.Columns.Item(2).ColumnWidth = 22.7109375
$ws.Columns.Item(3).ColumnWidth = 10.5703125
$ws.Columns.Item(4).ColumnWidth = 22.7109375
$ws.Columns.Item(5).ColumnWidth = 9.28515625
